# A new weekly price record was added to the Cilantro report. It is
# inserted as row 28 (dated 2021-10-15), pushing the existing rows
# 28-52 down to 29-53.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(28).Insert()

$ws.Range("A28").Value = 1
$ws.Range("B28").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C28").Value = 'Arica y Parinacota'
$ws.Range("D28").Value = 44484
$ws.Range("E28").Value = 15
$ws.Range("F28").Value = 100112040
$ws.Range("G28").Value = 'Cilantro'
$ws.Range("H28").Value = 'Sin especificar'
$ws.Range("I28").Value = 'Primera'
$ws.Range("J28").Value = 250
$ws.Range("K28").Value = 950
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 975
$ws.Range("N28").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O28").Value = 'Región de Arica y Parinacota'
$ws.Range("P28").Value = 488
$ws.Range("Q28").Value = 2
$ws.Range("R28").Value = 'Hortaliza'
